# Apply updated crypto price/volume figures (Sat Dec 16 10:54:13 UTC 2023 run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text that often looks numeric (e.g. "246.71" or
# "42.237.15"). Plain assignment lets Excel auto-convert numeric-looking text
# into a real number, silently dropping formatting such as trailing zeros
# ("42.60" -> 42.6) and changing the stored cell type away from text. A leading
# apostrophe forces Excel to keep the entry as text (and is stripped from the
# stored value), matching the source data.
function Set-TextCell($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

Set-TextCell 2 4 "42.237.15"
$ws.Cells.Item(2, 5).Value = "  -1.52%  "

Set-TextCell 3 4 "2.246.71"
$ws.Cells.Item(3, 5).Value = "  -1.44%  "

$ws.Cells.Item(4, 5).Value = "  -0.02%  "

Set-TextCell 5 4 "246.71"
$ws.Cells.Item(5, 5).Value = "  -1.44%  "

Set-TextCell 6 4 "0.622"
$ws.Cells.Item(6, 5).Value = "  -1.48%  "

Set-TextCell 7 4 "74.57"
$ws.Cells.Item(7, 5).Value = "  -5.76%  "

$ws.Cells.Item(8, 5).Value = "  +0.06%  "

Set-TextCell 9 4 "0.619"
$ws.Cells.Item(9, 5).Value = "  -3.99%  "

Set-TextCell 10 4 "42.60"
$ws.Cells.Item(10, 5).Value = "  +2.97%  "

Set-TextCell 11 4 "0.0949"
$ws.Cells.Item(11, 5).Value = "  -2.49%  "

Set-TextCell 12 4 "7.17"
$ws.Cells.Item(12, 5).Value = "  -2.63%  "

$ws.Cells.Item(13, 5).Value = "  -1.72%  "

Set-TextCell 14 4 "14.49"
$ws.Cells.Item(14, 5).Value = "  -4.85%  "

Set-TextCell 15 4 "0.853"
$ws.Cells.Item(15, 5).Value = "  -2.21%  "

Set-TextCell 16 4 "2.265.92"
$ws.Cells.Item(16, 5).Value = "  -0.50%  "

Set-TextCell 17 4 "42.126.54"
$ws.Cells.Item(17, 5).Value = "  -1.35%  "

Set-TextCell 18 4 "0.0000102"
$ws.Cells.Item(18, 5).Value = "  +1.51%  "

Set-TextCell 19 4 "72.39"
$ws.Cells.Item(19, 5).Value = "  +0.21%  "

$ws.Cells.Item(20, 5).Value = "  -1.72%  "

Set-TextCell 21 4 "2.22"
$ws.Cells.Item(21, 5).Value = "  +1.94%  "

Set-TextCell 22 4 "231.63"
$ws.Cells.Item(22, 5).Value = "  -1.09%  "

Set-TextCell 23 4 "8.90"
$ws.Cells.Item(23, 5).Value = "  +33.18%  "

$ws.Cells.Item(24, 5).Value = "  +0.03%  "

Set-TextCell 25 4 "11.43"
$ws.Cells.Item(25, 5).Value = "  +0.34%  "

$ws.Cells.Item(26, 5).Value = "  -4.36%  "

$ws.Cells.Item(27, 5).Value = "  -1.41%  "

$ws.Cells.Item(28, 5).Value = "  +3.65%  "

Set-TextCell 29 4 "168.78"
$ws.Cells.Item(29, 5).Value = "  -0.07%  "

Set-TextCell 30 4 "20.70"
$ws.Cells.Item(30, 5).Value = "  -0.90%  "

$ws.Cells.Item(31, 5).Value = "  -3.88%  "

$ws.Cells.Item(32, 5).Value = "  -1.08%  "

Set-TextCell 33 4 "30.84"
$ws.Cells.Item(33, 5).Value = "  +1.01%  "

$ws.Cells.Item(34, 5).Value = "  -1.75%  "

Set-TextCell 35 4 "5.25"
$ws.Cells.Item(35, 5).Value = "  +9.89%  "

$ws.Cells.Item(36, 5).Value = "  -2.71%  "

$ws.Cells.Item(37, 5).Value = "  +3.89%  "

Set-TextCell 38 4 "13.83"
$ws.Cells.Item(38, 5).Value = "  +1.79%  "

$ws.Cells.Item(39, 5).Value = "  -4.18%  "

Set-TextCell 40 4 "5.80"
$ws.Cells.Item(40, 5).Value = "  -2.43%  "

Set-TextCell 41 4 "62.69"
$ws.Cells.Item(41, 5).Value = "  +1.93%  "

$ws.Cells.Item(42, 5).Value = "  -3.19%  "

Set-TextCell 43 4 "106.90"
$ws.Cells.Item(43, 5).Value = "  -7.31%  "

$ws.Cells.Item(44, 5).Value = "  +1.19%  "

$ws.Cells.Item(45, 5).Value = "  -2.88%  "

$ws.Cells.Item(46, 5).Value = "  -0.41%  "

$ws.Cells.Item(47, 5).Value = "  +0.18%  "

$ws.Cells.Item(48, 5).Value = "  -3.64%  "

$ws.Cells.Item(49, 5).Value = "  +2.21%  "

Set-TextCell 50 4 "4.21"
$ws.Cells.Item(50, 5).Value = "  -10.56%  "

$ws.Cells.Item(51, 5).Value = "  -0.20%  "
